$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 2.24
$ws.Range("I3").Value = 55.31
$ws.Range("I4").Value = 15
$ws.Range("I5").Value = 6.72
$ws.Range("I6").Value = 120.48
$ws.Range("J6").Value = 114.24
$ws.Range("I7").Value = 31.35
$ws.Range("I8").Value = 60.46
